$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 already held the first "application mobile" entry; the author
# rewrote the note (it now talks about the QR-code scan feature instead of
# the old "requetes" blurb) and kept the same "6 heures" duration.
$ws.Range("B8").Value = "Commencement de l'application mobile. J'ai ajouté les champs qui apparaitront sur les deux pages. J'ai aussi ajouté la partie qui me permet de scanner un QR Code et que l'application me renvoie des informations que j'ai choisies"
$ws.Range("C8").Value = "6 heures"

# New day's journal entry (row 9): date, note and hours worked.
$ws.Range("A9").Value = 43236
$ws.Range("B9").Value = "J'ai commencé à essayer de faire des requêtes depuis l'application mobile vers ma base de données mais je n'y arrive pas encore.`n Il y a encore un problème avec les variables. M. Chavey est venu m'aider mais en vain "
$ws.Range("C9").Value = "5 heures"

# Match the date formatting already used by column A (re-use the style
# instead of typing a new number format so it keeps the same style id).
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Wrap text on the notes column, like the other multi-line entries.
$ws.Range("B9").WrapText = $true

# Taller row to fit the longer, two-line note.
$ws.Rows.Item(9).RowHeight = 45

# Column B widens slightly to accommodate the new text.
$ws.Columns.Item(2).ColumnWidth = 111.6

$ws.Range("C10").Select()
